# This script applies a rotation of row data among rows
# 10,12,13,14,15,16,17,18,19,20,21,22,24,25,26,27 on the active sheet,
# as described by the target XML diff: each destination row ends up
# holding the values that a different ("source") row used to hold.
#
# Only the columns whose contents can actually differ between these rows
# are touched (A,B,D,E,F,G,H,J,K,L,M,N,Q,R,AF); the remaining columns are
# identical across all the involved rows in the original workbook, so
# they are deliberately left untouched to avoid any unintended
# reinterpretation of their (text-looking-like-dates) values by Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Contiguous column blocks that need to be relocated.
$blocks = @("A:B", "D:H", "J:N", "Q:R", "AF:AF")

# The set of rows participating in the rotation.
$rows = @(10,12,13,14,15,16,17,18,19,20,21,22,24,25,26,27)

# destination row -> source row (i.e. destination row ends up holding the
# data that previously lived in the source row)
$mapping = @{
    10 = 16
    12 = 14
    13 = 19
    14 = 12
    15 = 13
    16 = 22
    17 = 18
    18 = 10
    19 = 20
    20 = 17
    21 = 15
    22 = 21
    24 = 26
    25 = 27
    26 = 24
    27 = 25
}

# 1) Snapshot the current ("before") contents of every involved row, for
#    each relevant column block.
$snapshot = @{}
foreach ($r in $rows) {
    $rowSnap = @{}
    foreach ($block in $blocks) {
        $parts = $block.Split(":")
        $addr = "$($parts[0])$r`:$($parts[1])$r"
        $rowSnap[$block] = $ws.Range($addr).Value2
    }
    $snapshot[$r] = $rowSnap
}

# 2) Write back each destination row using the snapshot taken from its
#    mapped source row, so every write uses pristine "before" data.
foreach ($destRow in $rows) {
    $srcRow = $mapping[$destRow]
    $srcSnap = $snapshot[$srcRow]
    foreach ($block in $blocks) {
        $parts = $block.Split(":")
        $addr = "$($parts[0])$destRow`:$($parts[1])$destRow"
        $ws.Range($addr).Value2 = $srcSnap[$block]
    }
}
